# Edit script implementing:
#   "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Net effect on the worksheet ("Hoja1"):
#   - A new arrears period "2509" is added as a new data row (after the
#     existing last row "2508"), reusing the bottom-border row style while
#     the former last row ("2508") is converted to a regular middle row.
#   - The "Periodo Mora" column (E) for all data rows gets centered
#     horizontal alignment.
#   - The summary values are refreshed: VALOR MORA total (E11) and
#     Cant. Periodos (F13) grow to account for the new period.
#   - The two closing/signature lines move down one row (27->28, 28->29)
#     to make room for the newly inserted data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a blank row at 23. This shifts the trailing "firma" rows
#     (old 27 -> 28, old 28 -> 29) down automatically, while row 22 (the
#     current last data row, "2508") is left untouched for now. ---
$ws.Rows(23).Insert()

# --- Duplicate the old last-row (row 22, still holding its original
#     bottom-border formatting) into the freshly inserted row 23, so the
#     new period keeps the "bottom of table" look. ---
$ws.Range("B22:J22").Copy($ws.Range("B23:J23"))

# --- Row 22 is no longer the last row, so restyle it like the other
#     interior rows (copy formatting only from row 21). ---
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 23 becomes the new "2509" period entry; other columns (worker
#     id/name, Valor Mora, Salario Basico) stay identical to the rest of
#     the table, and the Novedad/Observaciones cells stay blank. ---
$ws.Range("E23").Value = "2509"

# --- New formatting: center the "Periodo Mora" column for every data
#     row, including the newly added one. ---
$ws.Range("E16:E23").HorizontalAlignment = -4108

# --- Refresh header summary figures for the extra period. ---
$ws.Range("E11").Value = 290728
$ws.Range("F13").Value = 8
